$wb = $excel.ActiveWorkbook

# --- zh-cn sheet: update handoff/handback datetimes for file 2a79a01b... (row 2) ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("H2").Value = "2016-09-06 09:37:28"
$wsZh.Range("K2").Value = "2016-09-06 09:38:23"

# --- de-de sheet: update handoff/handback datetimes for file 2a79a01b... (row 2) ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("H2").Value = "2016-09-06 09:37:39"
$wsDe.Range("K2").Value = "2016-09-06 09:38:40"

# --- Overview sheet: "Latest HO Xliff Generate Date" reflects the de-de handoff datetime ---
$wsOv = $wb.Worksheets.Item("Overview")
$wsOv.Range("G2").Value = "2016-09-06 09:37:39"
$wsOv.Range("G3").Value = "2016-09-06 09:35:21"
